$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 2 and 3 had their date / volume / price values swapped.
# Row 2 (date 2021-06-10 / 44365) becomes row 2 with date 2021-04-28 / 44322 and row 3 values move the other way.

$ws.Range("D2").Value = 44322
$ws.Range("M2").Value = 600
$ws.Range("N2").Value = 1500
$ws.Range("O2").Value = 1600
$ws.Range("P2").Value = 1550
$ws.Range("S2").Value = 1550

$ws.Range("D3").Value = 44365
$ws.Range("M3").Value = 900
$ws.Range("N3").Value = 1200
$ws.Range("O3").Value = 1400
$ws.Range("P3").Value = 1300
$ws.Range("S3").Value = 1300
